# Update sheet name to reflect new "through" date
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2022-06-06"

# Update the June row label text (shared string) - cell A7
$ws.Range("A7").Value = "June (through 06-06)"

# Update June row (row 7) values for columns C (2016), G (2020), H (2021), I (2022)
$ws.Range("C7").Value = 8
$ws.Range("G7").Value = 36
$ws.Range("H7").Value = 26
$ws.Range("I7").Value = 18

# Update Total row (row 8) values for columns C (2016), G (2020), H (2021), I (2022)
$ws.Range("C8").Value = 217
$ws.Range("G8").Value = 394
$ws.Range("H8").Value = 657
$ws.Range("I8").Value = 682
